# Update the FOB price placeholder text on the slide:
#   "FOB : ${FOB}"  ->  "FOB : ${FOB Price}"

$p = $ppt.ActivePresentation

$oldText = 'FOB : ${FOB}'
$newText = 'FOB : ${FOB Price}'

foreach ($s in $p.Slides) {
    foreach ($shape in $s.Shapes) {
        if ($shape.HasTextFrame) {
            $tf = $shape.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                $full = $tr.Text
                $idx = $full.IndexOf($oldText)
                while ($idx -ge 0) {
                    # Narrow the edit to just the matched characters so the
                    # surrounding run/formatting is left untouched.
                    $sub = $tr.Characters($idx + 1, $oldText.Length)
                    $sub.Text = $newText

                    $full = $tr.Text
                    $idx = $full.IndexOf($oldText)
                }
            }
        }
    }
}
